# Auto-generated edit script updating Leve profit/price figures per scheduled runner refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 7815013.5
$ws.Range("I33").Value = 12500262
$ws.Range("J33").Value = 6266
$ws.Range("K33").Value = 12500262
$ws.Range("L33").Value = 6266
$ws.Range("M33").Value = -12500033
$ws.Range("N33").Value = -6724

$ws.Range("H40").Value = 614.4
$ws.Range("J40").Value = 614.4
$ws.Range("L40").Value = 614.4
$ws.Range("N40").Value = -964.4

$ws.Range("H70").Value = 1318.7142
$ws.Range("J70").Value = 1483
$ws.Range("L70").Value = 4449
$ws.Range("N70").Value = -4989

$ws.Range("H73").Value = 1318.7142
$ws.Range("J73").Value = 1483
$ws.Range("L73").Value = 4449
$ws.Range("N73").Value = -6321

$ws.Range("H74").Value = 2876.1875
$ws.Range("I74").Value = 2219.2144
$ws.Range("J74").Value = 7475
$ws.Range("K74").Value = 2219.2144
$ws.Range("L74").Value = 7475
$ws.Range("M74").Value = -1283.2144
$ws.Range("N74").Value = -9347

$ws.Range("H77").Value = 2876.1875
$ws.Range("I77").Value = 2219.2144
$ws.Range("J77").Value = 7475
$ws.Range("K77").Value = 11096.072
$ws.Range("L77").Value = 37375
$ws.Range("M77").Value = -6416.072
$ws.Range("N77").Value = -46735

$ws.Range("H103").Value = 1385.9
$ws.Range("J103").Value = 1563.875
$ws.Range("L103").Value = 4691.625
$ws.Range("N103").Value = -5863.625

$ws.Range("H112").Value = 3984.0527
$ws.Range("I112").Value = 8255.556
$ws.Range("K112").Value = 24766.668
$ws.Range("M112").Value = -23658.668

$ws.Range("H132").Value = 5386.44
$ws.Range("I132").Value = 3638.8096
$ws.Range("K132").Value = 10916.4288
$ws.Range("M132").Value = -8386.4288

$ws.Range("H137").Value = 11590.546
$ws.Range("J137").Value = 18334.652
$ws.Range("L137").Value = 55003.95599999999
$ws.Range("N137").Value = -60103.95599999999

$ws.Range("H138").Value = 4262.6743
$ws.Range("I138").Value = 3302.7
$ws.Range("J138").Value = 4553.5757
$ws.Range("K138").Value = 9908.099999999999
$ws.Range("L138").Value = 13660.7271
$ws.Range("M138").Value = -4768.099999999999
$ws.Range("N138").Value = -23940.7271

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1149.7142
$ws.Range("I102").Value = 1149.7142
$ws.Range("K102").Value = 1149.7142
$ws.Range("M102").Value = 472.2858000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1144.5333
$ws.Range("I105").Value = 937.1923
$ws.Range("J105").Value = 2492.25
$ws.Range("K105").Value = 937.1923
$ws.Range("L105").Value = 2492.25
$ws.Range("M105").Value = 809.8077
$ws.Range("N105").Value = -5986.25

$ws.Range("H107").Value = 923.25
$ws.Range("I107").Value = 899.3333
$ws.Range("J107").Value = 995
$ws.Range("K107").Value = 899.3333
$ws.Range("L107").Value = 995
$ws.Range("M107").Value = 1020.6667
$ws.Range("N107").Value = -4835

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1687.9
$ws.Range("I22").Value = 776.8
$ws.Range("K22").Value = 776.8
$ws.Range("M22").Value = -426.8

$ws.Range("H93").Value = 9999
$ws.Range("I93").Value = 9999
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 9999
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -8127
$ws.Range("N93").Value = ""

$ws.Range("H105").Value = 7536.8
$ws.Range("I105").Value = 8003.7144
$ws.Range("K105").Value = 8003.7144
$ws.Range("M105").Value = -6256.7144

$ws.Range("H107").Value = 3198.2173
$ws.Range("I107").Value = 1238.0588
$ws.Range("K107").Value = 1238.0588
$ws.Range("M107").Value = 681.9412

$ws.Range("H134").Value = 20004180
$ws.Range("I134").Value = 1288.8948
$ws.Range("J134").Value = 83346664
$ws.Range("K134").Value = 3866.6844
$ws.Range("L134").Value = 250039992
$ws.Range("M134").Value = -1331.6844
$ws.Range("N134").Value = -250045062

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1081.5714
$ws.Range("I5").Value = 1220.75
$ws.Range("J5").Value = 896
$ws.Range("K5").Value = 3662.25
$ws.Range("L5").Value = 2688
$ws.Range("M5").Value = -3550.25
$ws.Range("N5").Value = -2912

$ws.Range("H31").Value = 4829.6
$ws.Range("J31").Value = 9999
$ws.Range("L31").Value = 29997
$ws.Range("N31").Value = -30573

$ws.Range("H68").Value = 2536.4546
$ws.Range("J68").Value = 2562
$ws.Range("L68").Value = 7686
$ws.Range("N68").Value = -9308

$ws.Range("H71").Value = 2536.4546
$ws.Range("J71").Value = 2562
$ws.Range("L71").Value = 23058
$ws.Range("N71").Value = -31170

$ws.Range("H80").Value = 27143.46
$ws.Range("I80").Value = 19123
$ws.Range("K80").Value = 57369
$ws.Range("M80").Value = -56433

$ws.Range("H83").Value = 27143.46
$ws.Range("I83").Value = 19123
$ws.Range("K83").Value = 172107
$ws.Range("M83").Value = -167427

$ws.Range("H107").Value = 11080.5
$ws.Range("I107").Value = 576.5
$ws.Range("J107").Value = 18083.166
$ws.Range("K107").Value = 1729.5
$ws.Range("L107").Value = 54249.49800000001
$ws.Range("M107").Value = 190.5
$ws.Range("N107").Value = -58089.49800000001

$ws.Range("H110").Value = 35000
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 35000
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 105000
$ws.Range("M110").Value = ""
$ws.Range("N110").Value = -113180

$ws.Range("H111").Value = 1299.5
$ws.Range("I111").Value = 1299.5
$ws.Range("K111").Value = 3898.5
$ws.Range("M111").Value = -831.5

$ws.Range("H132").Value = 2332729.2
$ws.Range("I132").Value = 1593.75
$ws.Range("K132").Value = 14343.75
$ws.Range("M132").Value = -11813.75

$ws.Range("H134").Value = 3297.5925
$ws.Range("I134").Value = 1528.1578
$ws.Range("K134").Value = 4584.4734
$ws.Range("M134").Value = 485.5266000000001

$ws.Range("H135").Value = 1081.5714
$ws.Range("I135").Value = 1220.75
$ws.Range("J135").Value = 896
$ws.Range("K135").Value = 10986.75
$ws.Range("L135").Value = 8064
$ws.Range("M135").Value = -8451.75
$ws.Range("N135").Value = -13134

$ws.Range("H140").Value = 2110.2222
$ws.Range("I140").Value = 1899
$ws.Range("K140").Value = 5697
$ws.Range("M140").Value = -517

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 130.63637
$ws.Range("I2").Value = 161
$ws.Range("J2").Value = 49.666668
$ws.Range("K2").Value = 161
$ws.Range("L2").Value = 49.666668
$ws.Range("M2").Value = -48
$ws.Range("N2").Value = -275.666668

$ws.Range("H126").Value = 18206.143
$ws.Range("I126").Value = 13888.8
$ws.Range("K126").Value = 41666.39999999999
$ws.Range("M126").Value = -39196.39999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4197.8125
$ws.Range("I22").Value = 685.2
$ws.Range("J22").Value = 5794.4546
$ws.Range("K22").Value = 685.2
$ws.Range("L22").Value = 5794.4546
$ws.Range("M22").Value = -390.2
$ws.Range("N22").Value = -6384.4546

$ws.Range("H27").Value = 4197.8125
$ws.Range("I27").Value = 685.2
$ws.Range("J27").Value = 5794.4546
$ws.Range("K27").Value = 685.2
$ws.Range("L27").Value = 5794.4546
$ws.Range("M27").Value = -578.2
$ws.Range("N27").Value = -6008.4546

$ws.Range("H55").Value = 162.7619
$ws.Range("I55").Value = 227.375
$ws.Range("J55").Value = 123
$ws.Range("K55").Value = 227.375
$ws.Range("L55").Value = 123
$ws.Range("M55").Value = -54.375
$ws.Range("N55").Value = -469

$ws.Range("H56").Value = 10042.5
$ws.Range("I56").Value = 12000
$ws.Range("K56").Value = 12000
$ws.Range("M56").Value = -11309

$ws.Range("H100").Value = 6804.4546
$ws.Range("I100").Value = 6262.5
$ws.Range("K100").Value = 6262.5
$ws.Range("M100").Value = -5721.5

$ws.Range("H122").Value = 7740.452
$ws.Range("I122").Value = 6516.04
$ws.Range("K122").Value = 19548.12
$ws.Range("M122").Value = -17098.12

$ws.Range("H136").Value = 1286659.8
$ws.Range("I136").Value = 27833.1
$ws.Range("J136").Value = 2125877.5
$ws.Range("K136").Value = 83499.29999999999
$ws.Range("L136").Value = 6377632.5
$ws.Range("M136").Value = -80949.29999999999
$ws.Range("N136").Value = -6382732.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2094.8333
$ws.Range("I81").Value = 2094.8333
$ws.Range("K81").Value = 4189.6666
$ws.Range("M81").Value = -3128.6666

$ws.Range("H84").Value = 2094.8333
$ws.Range("I84").Value = 2094.8333
$ws.Range("K84").Value = 20948.333
$ws.Range("M84").Value = -15644.333

$ws.Range("H107").Value = 10768.2
$ws.Range("I107").Value = 999
$ws.Range("K107").Value = 2997
$ws.Range("M107").Value = -1077

$ws.Range("H125").Value = 50000
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 50000
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 50000
$ws.Range("M125").Value = ""
$ws.Range("N125").Value = -59840

$ws.Range("H136").Value = 327992.5
$ws.Range("I136").Value = 2675.2727
$ws.Range("K136").Value = 8025.8181
$ws.Range("M136").Value = -5475.8181

Write-Host "Updated 238 cells across 8 sheets"